# Adds a new "HISTORY" column (M) and records a bulk leave-selection update
# for a handful of employees (rows 3-11): their SOLDE_Y-1 (I), SOLDE_Y (J)
# and SOLDE (K) balances change, and a HISTORY string capturing the
# selected date range is written to the new column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("M1").Value = "HISTORY"

# Row 3
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 33
$ws.Range("M3").Value = "2025-07-19_2025-07-24"

# Row 4
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 28
$ws.Range("M4").Value = "2025-07-19_2025-07-24"

# Row 5
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 31
$ws.Range("M5").Value = "2025-07-19_2025-07-24"

# Row 6
$ws.Range("J6").Value = 24
$ws.Range("K6").Value = 24
$ws.Range("M6").Value = "2025-07-21_2025-07-27"

# Row 7
$ws.Range("J7").Value = 24
$ws.Range("K7").Value = 24
$ws.Range("M7").Value = "2025-07-21_2025-07-27"

# Row 8
$ws.Range("J8").Value = 24
$ws.Range("K8").Value = 24
$ws.Range("M8").Value = "2025-07-21_2025-07-27"

# Row 10
$ws.Range("J10").Value = 27
$ws.Range("K10").Value = 27
$ws.Range("M10").Value = "2025-07-25_2025-07-28"

# Row 11
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = "2025-08-01_2025-08-31"
